# Update odds figures in Sheet1 to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - Flamengo RJ vs Fluminense
$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25

# Row 4 - America MG vs Goias
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5

# Row 5 - Shrewsbury vs Exeter
$ws.Range("I5").Value = 2.25
$ws.Range("J5").Value = 3.75
$ws.Range("AA5").Value = 26
$ws.Range("AH5").Value = 7.5
$ws.Range("AJ5").Value = 9

# Row 8 - Union Comercio vs Comerciantes Unidos
$ws.Range("G8").Value = 2.2
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 2.88
$ws.Range("K8").Value = 2.1
$ws.Range("L8").Value = 3.75
$ws.Range("W8").Value = 8
$ws.Range("X8").Value = 11
$ws.Range("Y8").Value = 9.5
$ws.Range("Z8").Value = 21
$ws.Range("AH8").Value = 11
$ws.Range("AI8").Value = 17
$ws.Range("AJ8").Value = 12
$ws.Range("AK8").Value = 34
$ws.Range("AL8").Value = 26
$ws.Range("AM8").Value = 34
$ws.Range("AN8").Value = 4.33
$ws.Range("AO8").Value = 12
$ws.Range("AW8").Value = 5
$ws.Range("AX8").Value = 19
$ws.Range("AY8").Value = 26
$ws.Range("AZ8").Value = 51

# Row 9 - Cusco vs Los Chankas
$ws.Range("M9").Value = 1.02
$ws.Range("N9").Value = 19
$ws.Range("Q9").Value = 1.5
$ws.Range("R9").Value = 2.5

# Row 10 - Sport Huancayo vs Grau
$ws.Range("G10").Value = 2.05
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 3.7
$ws.Range("Q10").Value = 2.3
$ws.Range("R10").Value = 1.6
$ws.Range("W10").Value = 6
$ws.Range("AU10").Value = 9
$ws.Range("AZ10").Value = 81
$ws.Range("BD10").Value = 151

# Row 12 - Fenix vs CA Cerro
$ws.Range("K12").Value = 1.95
$ws.Range("AB12").Value = 41
$ws.Range("AE12").Value = 21
$ws.Range("AF12").Value = 81
$ws.Range("AX12").Value = 23
$ws.Range("BB12").Value = 351

# Row 13 - Nacional vs Miramar
$ws.Range("G13").Value = 1.29
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 12
$ws.Range("J13").Value = 1.73
$ws.Range("K13").Value = 2.63
$ws.Range("L13").Value = 9
$ws.Range("O13").Value = 1.18
$ws.Range("P13").Value = 4.5
$ws.Range("Q13").Value = 1.62
$ws.Range("R13").Value = 2.25
$ws.Range("U13").Value = 2.1
$ws.Range("V13").Value = 1.67
$ws.Range("X13").Value = 6.5
$ws.Range("AB13").Value = 29
$ws.Range("AF13").Value = 67
$ws.Range("AS13").Value = 126
$ws.Range("AW13").Value = 10
$ws.Range("AX13").Value = 41
$ws.Range("AY13").Value = 41
$ws.Range("AZ13").Value = 201
$ws.Range("BA13").Value = 201
